# Custom Employee Template Export
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table data (4 rows x 8 cols) replacing the previous 6-row table
$data = @(
    @(58, "Ananda Zakia Syahfitri", 1, 12345, 1, "Sekretaris", "2023-12-31T04:45:40.000000Z", "2023-12-31T04:45:40.000000Z"),
    @(59, "Bima Sakti", 4, 6789, 3, "Anggota", "2023-12-31T04:55:18.000000Z", "2023-12-31T04:55:18.000000Z"),
    @(60, "Ananda Zakia Syahfitri", 1, 12345, 1, "Sekretaris", "2023-12-31T04:56:47.000000Z", "2023-12-31T04:56:47.000000Z"),
    @(61, "Bima Sakti", 4, 6789, 3, "Anggota", "2023-12-31T04:56:47.000000Z", "2023-12-31T04:56:47.000000Z")
)

# Clear out old rows (rows 5 and 6 no longer exist)
$ws.Range("A5:H6").Clear()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
